$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 5480.087
$ws.Range("I86").Value = 1144.6666
$ws.Range("J86").Value = 10209.637
$ws.Range("K86").Value = 1144.6666
$ws.Range("L86").Value = 10209.637
$ws.Range("M86").Value = -21.66660000000002
$ws.Range("N86").Value = -12455.637
# Row 89
$ws.Range("H89").Value = 5480.087
$ws.Range("I89").Value = 1144.6666
$ws.Range("J89").Value = 10209.637
$ws.Range("K89").Value = 5723.333000000001
$ws.Range("L89").Value = 51048.185
$ws.Range("M89").Value = -107.3330000000005
$ws.Range("N89").Value = -62280.185
# Row 95
$ws.Range("H95").Value = 36500
$ws.Range("J95").Value = 36500
$ws.Range("L95").Value = 36500
$ws.Range("N95").Value = -41992
# Row 121
$ws.Range("H121").Value = 6362.7617
$ws.Range("I121").Value = 275
$ws.Range("J121").Value = 7003.579
$ws.Range("K121").Value = 825
$ws.Range("L121").Value = 21010.737
$ws.Range("M121").Value = 922
$ws.Range("N121").Value = -24504.737
# Row 129
$ws.Range("H129").Value = 157233.44
$ws.Range("J129").Value = 176506.34
$ws.Range("L129").Value = 529519.02
$ws.Range("N129").Value = -539519.02

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1788.2222
$ws.Range("I2").Value = 1882.2
$ws.Range("J2").Value = 1670.75
$ws.Range("K2").Value = 1882.2
$ws.Range("L2").Value = 1670.75
$ws.Range("M2").Value = -1769.2
$ws.Range("N2").Value = -1896.75
# Row 4
$ws.Range("H4").Value = 189.23077
$ws.Range("I4").Value = 132.5
$ws.Range("J4").Value = 280
$ws.Range("K4").Value = 132.5
$ws.Range("L4").Value = 280
$ws.Range("M4").Value = -16.5
$ws.Range("N4").Value = -512
# Row 32
$ws.Range("H32").Value = 5499.864
$ws.Range("I32").Value = 4140.8228
$ws.Range("K32").Value = 4140.8228
$ws.Range("M32").Value = -3853.8228
# Row 45
$ws.Range("H45").Value = 2830.55
$ws.Range("I45").Value = 3012.3
$ws.Range("K45").Value = 3012.3
$ws.Range("M45").Value = -2635.3
# Row 74
$ws.Range("H74").Value = 30304552
$ws.Range("I74").Value = 40000630
$ws.Range("K74").Value = 40000630
$ws.Range("M74").Value = -39999756
# Row 77
$ws.Range("H77").Value = 30304552
$ws.Range("I77").Value = 40000630
$ws.Range("K77").Value = 200003150
$ws.Range("M77").Value = -199998782
# Row 116
$ws.Range("H116").Value = 1788.2222
$ws.Range("I116").Value = 1882.2
$ws.Range("J116").Value = 1670.75
$ws.Range("K116").Value = 1882.2
$ws.Range("L116").Value = 1670.75
$ws.Range("M116").Value = 411.8
$ws.Range("N116").Value = -6258.75
# Row 132
$ws.Range("H132").Value = 26284.637
$ws.Range("I132").Value = 3941.077
$ws.Range("J132").Value = 58558.668
$ws.Range("K132").Value = 11823.231
$ws.Range("L132").Value = 175676.004
$ws.Range("M132").Value = -9293.231
$ws.Range("N132").Value = -180736.004
# Row 139
$ws.Range("H139").Value = 38794.832
$ws.Range("J139").Value = 38794.832
$ws.Range("L139").Value = 38794.832
$ws.Range("N139").Value = -49074.832

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1788.2222
$ws.Range("I3").Value = 1882.2
$ws.Range("J3").Value = 1670.75
$ws.Range("K3").Value = 1882.2
$ws.Range("L3").Value = 1670.75
$ws.Range("M3").Value = -1768.2
$ws.Range("N3").Value = -1898.75
# Row 20
$ws.Range("H20").Value = 1504.2683
$ws.Range("I20").Value = 1432.2609
$ws.Range("J20").Value = 1596.2778
$ws.Range("K20").Value = 1432.2609
$ws.Range("L20").Value = 1596.2778
$ws.Range("M20").Value = -1185.2609
$ws.Range("N20").Value = -2090.2778

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 858
$ws.Range("I107").Value = 322.22726
$ws.Range("K107").Value = 322.22726
$ws.Range("M107").Value = 1597.77274
# Row 132
$ws.Range("H132").Value = 7459.4287
$ws.Range("I132").Value = 5720.4
$ws.Range("J132").Value = 11807
$ws.Range("K132").Value = 17161.2
$ws.Range("L132").Value = 35421
$ws.Range("M132").Value = -14631.2
$ws.Range("N132").Value = -40481
# Row 134
$ws.Range("H134").Value = 1108.4
$ws.Range("I134").Value = 885.5
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2656.5
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -121.5
$ws.Range("N134").Value = -11070
# Row 137
$ws.Range("H137").Value = 20876
$ws.Range("J137").Value = 23845
$ws.Range("L137").Value = 23845
$ws.Range("N137").Value = -34045

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 398.5
$ws.Range("I7").Value = 45.5
$ws.Range("K7").Value = 136.5
$ws.Range("M7").Value = -24.5
# Row 39
$ws.Range("H39").Value = 4625
$ws.Range("J39").Value = 4625
$ws.Range("L39").Value = 13875
$ws.Range("N39").Value = -14463
# Row 69
$ws.Range("H69").Value = 2116.1667
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 2424.25
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 7272.75
$ws.Range("M69").Value = -3689
$ws.Range("N69").Value = -8894.75
# Row 72
$ws.Range("H72").Value = 2116.1667
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 2424.25
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 21818.25
$ws.Range("M72").Value = -9444
$ws.Range("N72").Value = -29930.25
# Row 122
$ws.Range("H122").Value = 912.5
$ws.Range("J122").Value = 1005.1429
$ws.Range("L122").Value = 9046.286100000001
$ws.Range("N122").Value = -13946.2861
# Row 131
$ws.Range("H131").Value = 705.84
$ws.Range("I131").Value = 470
$ws.Range("J131").Value = 708.2222
$ws.Range("K131").Value = 1410
$ws.Range("L131").Value = 2124.6666
$ws.Range("M131").Value = 3630
$ws.Range("N131").Value = -12204.6666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3314.04
$ws.Range("I80").Value = 2938.9285
$ws.Range("J80").Value = 3791.4546
$ws.Range("K80").Value = 2938.9285
$ws.Range("L80").Value = 3791.4546
$ws.Range("M80").Value = -1940.9285
$ws.Range("N80").Value = -5787.4546
# Row 83
$ws.Range("H83").Value = 3314.04
$ws.Range("I83").Value = 2938.9285
$ws.Range("J83").Value = 3791.4546
$ws.Range("K83").Value = 14694.6425
$ws.Range("L83").Value = 18957.273
$ws.Range("M83").Value = -9702.6425
$ws.Range("N83").Value = -28941.273
# Row 97
$ws.Range("H97").Value = 1212.1143
$ws.Range("I97").Value = 1301.5927
$ws.Range("J97").Value = 910.125
$ws.Range("K97").Value = 1301.5927
$ws.Range("L97").Value = 910.125
$ws.Range("M97").Value = -805.5926999999999
$ws.Range("N97").Value = -1902.125
# Row 132
$ws.Range("H132").Value = 65633.664
$ws.Range("I132").Value = 22008
$ws.Range("J132").Value = 87446.5
$ws.Range("K132").Value = 66024
$ws.Range("L132").Value = 262339.5
$ws.Range("M132").Value = -63494
$ws.Range("N132").Value = -267399.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5329.353
$ws.Range("I7").Value = 5230.6924
$ws.Range("J7").Value = 5650
$ws.Range("K7").Value = 5230.6924
$ws.Range("L7").Value = 5650
$ws.Range("M7").Value = -5118.6924
$ws.Range("N7").Value = -5874
# Row 22
$ws.Range("H22").Value = 3408.5557
$ws.Range("I22").Value = 4804.3335
$ws.Range("J22").Value = 617
$ws.Range("K22").Value = 4804.3335
$ws.Range("L22").Value = 617
$ws.Range("M22").Value = -4509.3335
$ws.Range("N22").Value = -1207
# Row 27
$ws.Range("H27").Value = 3408.5557
$ws.Range("I27").Value = 4804.3335
$ws.Range("J27").Value = 617
$ws.Range("K27").Value = 4804.3335
$ws.Range("L27").Value = 617
$ws.Range("M27").Value = -4697.3335
$ws.Range("N27").Value = -831
# Row 106
$ws.Range("H106").Value = 22185
$ws.Range("J106").Value = 22185
$ws.Range("L106").Value = 22185
$ws.Range("N106").Value = -24709
# Row 126
$ws.Range("H126").Value = 5329.353
$ws.Range("I126").Value = 5230.6924
$ws.Range("J126").Value = 5650
$ws.Range("K126").Value = 15692.0772
$ws.Range("L126").Value = 16950
$ws.Range("M126").Value = -13222.0772
$ws.Range("N126").Value = -21890
# Row 132
$ws.Range("H132").Value = 1923.2433
$ws.Range("I132").Value = 1472.9032
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 4418.7096
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -1888.7096
$ws.Range("N132").Value = -17810

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16144
# Row 92
$ws.Range("H92").Value = 30500
$ws.Range("J92").Value = 30500
$ws.Range("L92").Value = 30500
$ws.Range("N92").Value = -35492
# Row 132
$ws.Range("H132").Value = 757.8
$ws.Range("I132").Value = 578.2917
$ws.Range("J132").Value = 1475.8334
$ws.Range("K132").Value = 1734.8751
$ws.Range("L132").Value = 4427.5002
$ws.Range("M132").Value = 795.1249
$ws.Range("N132").Value = -9487.5002
# Row 136
$ws.Range("H136").Value = 52637030
$ws.Range("J136").Value = 12250
$ws.Range("L136").Value = 36750
$ws.Range("M136").Value = -41850
$ws.Range("N136").Value = -41850
